# Update cryptocurrency price/volume data (commit: "Updated cryptos list")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.992.85"
$ws.Range("E2").Value = "'  +0.42%  "
$ws.Range("D3").Value = "'2.558.37"
$ws.Range("E3").Value = "'  +0.14%  "
$ws.Range("E4").Value = "'  +0.24%  "
$ws.Range("D5").Value = "'582.19"
$ws.Range("E5").Value = "'  +0.65%  "
$ws.Range("D6").Value = "'170.71"
$ws.Range("E6").Value = "'  +0.10%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E8").Value = "'  +1.88%  "
$ws.Range("D9").Value = "'2.558.96"
$ws.Range("E9").Value = "'  +0.22%  "
$ws.Range("E10").Value = "'  +0.00%  "
$ws.Range("E11").Value = "'  -0.36%  "
$ws.Range("D12").Value = "'0.361"
$ws.Range("E12").Value = "'  +2.95%  "
$ws.Range("D13").Value = "'4.94"
$ws.Range("E13").Value = "'  +2.09%  "
$ws.Range("D14").Value = "'3.029.71"
$ws.Range("E14").Value = "'  -0.82%  "
$ws.Range("D15").Value = "'71.136.74"
$ws.Range("E15").Value = "'  +0.87%  "
$ws.Range("E16").Value = "'  -1.43%  "
$ws.Range("D17").Value = "'25.48"
$ws.Range("E17").Value = "'  +1.03%  "
$ws.Range("D18").Value = "'2.555.51"
$ws.Range("E18").Value = "'  +0.67%  "
$ws.Range("D19").Value = "'11.62"
$ws.Range("E19").Value = "'  -0.98%  "
$ws.Range("E20").Value = "'  +3.79%  "
$ws.Range("D21").Value = "'356.59"
$ws.Range("E21").Value = "'  -1.88%  "
$ws.Range("D22").Value = "'3.96"
$ws.Range("E22").Value = "'  -1.30%  "
$ws.Range("D23").Value = "'2.05"
$ws.Range("E23").Value = "'  +2.85%  "
$ws.Range("E24").Value = "'  +0.29%  "
$ws.Range("D25").Value = "'70.63"
$ws.Range("E25").Value = "'  +0.87%  "
$ws.Range("D26").Value = "'4.08"
$ws.Range("E26").Value = "'  -2.04%  "
$ws.Range("D27").Value = "'9.09"
$ws.Range("E27").Value = "'  -2.07%  "
$ws.Range("E28").Value = "'  +0.90%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "'  -0.31%  "
$ws.Range("D30").Value = "'0.0₃0923"
$ws.Range("E30").Value = "'  -0.96%  "
$ws.Range("D31").Value = "'7.99"
$ws.Range("E31").Value = "'  +2.40%  "
$ws.Range("D32").Value = "'472.17"
$ws.Range("E32").Value = "'  -2.09%  "
$ws.Range("E33").Value = "'  -2.05%  "
$ws.Range("E34").Value = "'  -0.09%  "
$ws.Range("E35").Value = "'  +0.12%  "
$ws.Range("E36").Value = "'  +4.40%  "
$ws.Range("D37").Value = "'158.03"
$ws.Range("E37").Value = "'  +0.59%  "
$ws.Range("E38").Value = "'  +0.52%  "
$ws.Range("D39").Value = "'19.13"
$ws.Range("E39").Value = "'  +1.54%  "
$ws.Range("E40").Value = "'  -0.03%  "
$ws.Range("D41").Value = "'4.89"
$ws.Range("E41").Value = "'  +2.16%  "
$ws.Range("E42").Value = "'  +0.38%  "
$ws.Range("E43").Value = "'  -3.74%  "
$ws.Range("D44").Value = "'2.36"
$ws.Range("E44").Value = "'  -4.65%  "
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "'  -12.15%  "
$ws.Range("E46").Value = "'  +0.67%  "
$ws.Range("D47").Value = "'145.17"
$ws.Range("E47").Value = "'  -0.92%  "
$ws.Range("D48").Value = "'0.539"
$ws.Range("E48").Value = "'  +1.24%  "
$ws.Range("D49").Value = "'3.56"
$ws.Range("E49").Value = "'  -0.26%  "
$ws.Range("E50").Value = "'  -1.05%  "
$ws.Range("D51").Value = "'0.0740"
$ws.Range("E51").Value = "'  +1.13%  "
